$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Insert a new column before column B (shifts dbExcel/WebExcel columns right,
# and preserves the precise stored widths of the other columns untouched)
$ws.Columns("B").Insert()

# New header + value for the inserted "StatQuery" column
$ws.Range("B1").Value = "StatQuery"
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN ['Carcinosarcoma of the uterus'] OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"

# Match the wrap-text style used by A2
$ws.Range("B2").WrapText = $true

# New column should be as wide as column A (its neighbour)
$ws.Columns("B").ColumnWidth = $ws.Columns("A").ColumnWidth

# Update selection to A2
$ws.Range("A2").Select() | Out-Null
